# Adds two new "Mac-Address" user rows (Jane Smith, John Doe) to the bottom
# of the user_detail_h master data sheet.
#
# The last existing data row is row 30 (id 110029 / Carolyn Rodriguez).
# New rows 31 and 32 must end up with the same per-column formatting as the
# existing data rows (column D uses style index 2, column I uses style
# index 1). Plain "set .Value on a new row" does not carry that formatting
# in this engine, but inserting a row (like pressing the Excel "Insert"
# command on a selected row) duplicates the formatting of the surrounding
# row automatically - so we insert two rows above the last row, move the
# original last row's data into the first of those (restoring row 30 to its
# original content/format), and then populate the two freshly-formatted
# rows that follow (31 and 32) with the new records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 30

# Insert two new formatted rows directly above the current last row. This
# pushes the current row 30 down to row 32, and rows 30-31 become blank
# rows that inherit row 30's original formatting.
$ws.Range("A" + $lastRow + ":A" + ($lastRow + 1)).EntireRow().Insert()

# Copy the original last row's values (now sitting at row 32) back up onto
# row 30 - values only, so row 30 keeps the formatting it already has.
$ws.Range("A" + ($lastRow + 2) + ":L" + ($lastRow + 2)).Copy()
$ws.Range("A" + $lastRow + ":L" + $lastRow).PasteSpecial(-4163)
$excel.CutCopyMode = 0

# Clear out the now-duplicated old data from row 32, leaving it blank but
# still correctly formatted, ready to receive the second new record.
$ws.Range("A" + ($lastRow + 2) + ":L" + ($lastRow + 2)).ClearContents()

# Populate row 32 (John Doe) first so its strings are registered in the
# shared-string table before Jane Smith's.
$r2 = $lastRow + 2
$ws.Cells.Item($r2, 1).Value = 110031
$ws.Cells.Item($r2, 2).Value = 9317596767
$ws.Cells.Item($r2, 3).Value = "John Doe"
$ws.Cells.Item($r2, 4).Value = "john.doe@xyz.com"
$ws.Cells.Item($r2, 5).Value = 818876431
$ws.Cells.Item($r2, 6).Value = "ACT"
$ws.Cells.Item($r2, 7).Value = "eng"
$ws.Cells.Item($r2, 8).Value = "PWD"
$ws.Cells.Item($r2, 9).Value = $true
$ws.Cells.Item($r2, 10).Value = "superadmin"
$ws.Cells.Item($r2, 11).Value = "now()"
$ws.Cells.Item($r2, 12).Value = "now()"

# Populate row 31 (Jane Smith) next.
$r1 = $lastRow + 1
$ws.Cells.Item($r1, 1).Value = 110030
$ws.Cells.Item($r1, 2).Value = 9317596768
$ws.Cells.Item($r1, 3).Value = "Jane Smith"
$ws.Cells.Item($r1, 4).Value = "jane.smith@xyz.com"
$ws.Cells.Item($r1, 5).Value = 818876432
$ws.Cells.Item($r1, 6).Value = "ACT"
$ws.Cells.Item($r1, 7).Value = "eng"
$ws.Cells.Item($r1, 8).Value = "PWD"
$ws.Cells.Item($r1, 9).Value = $true
$ws.Cells.Item($r1, 10).Value = "superadmin"
$ws.Cells.Item($r1, 11).Value = "now()"
$ws.Cells.Item($r1, 12).Value = "now()"

# Match the author's final cursor position.
$ws.Range("F30").Select() | Out-Null
